# Update the alt-text (OOXML descr attribute) of the chart picture shapes
# on several slides. Each picture is named "Picture 3" and is shape index 3
# on its respective slide.

$p = $ppt.ActivePresentation

# Slide 3: sales chart picture
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(3)
$sh.AlternativeText = "slide_89f29345_create_sales_chart.png"

# Slide 4: market share picture
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(3)
$sh.AlternativeText = "slide_efe173e6_create_market_share.png"

# Slide 11: growth trend picture
$s = $p.Slides.Item(11)
$sh = $s.Shapes.Item(3)
$sh.AlternativeText = "slide_27b543ea_create_growth_trend.png"

# Slide 13: sales chart picture (second occurrence)
$s = $p.Slides.Item(13)
$sh = $s.Shapes.Item(3)
$sh.AlternativeText = "slide_89f29345_create_sales_chart.png"
